# Code Version 190520 15:01
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Header row: drop the trailing "?" from the IsFolder/IsExist column headers
$ws.Range("D1").Value = "IsFolder"
$ws.Range("E1").Value = "IsExist"

# Fix the TSYSDW path: correct the leading directory and the doubled slash
$ws.Range("C2").Value = "/tsys/prime/UBP_ubp/datawarehouse_file/EDS_TRANSFER/TSYS_DW.sh"
